$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear D2, D3 values (were 0, now empty)
$ws.Range("D2:D3").ClearContents()

# Clear C4:C11 values (were 0, now empty)
$ws.Range("C4:C11").ClearContents()

# Update selection to D2:D3, active cell D2
$ws.Range("D2:D3").Select()
